$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.730.69'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.074.23'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.80%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.623'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.55'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.31%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +1.00%  '
$ws.Range('E10').Value = '  +0.67%  '
$ws.Range('E11').Value = '  +3.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.381.14'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.74'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.95'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.774'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.37'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.050.89'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.69%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.658.32'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.12'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.07'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0835'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.22'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('E25').Value = '  -2.61%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '170.89'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.94%  '
$ws.Range('E27').Value = '  +2.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.03'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.49'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('E30').Value = '  -2.42%  '
$ws.Range('E31').Value = '  +2.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.69'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.49%  '
$ws.Range('E33').Value = '  +1.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.65'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.18%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.49'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.74%  '
$ws.Range('E36').Value = '  +0.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.39'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.32'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '100.20'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0975'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.43%  '
$ws.Range('E42').Value = '  -1.93%  '
$ws.Range('E43').Value = '  +0.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.73'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +8.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.441.68'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.43%  '
$ws.Range('E46').Value = '  -0.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.21'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.59%  '
$ws.Range('E48').Value = '  +0.29%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.43'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.09%  '
$ws.Range('E50').Value = '  -1.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.265.89'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.84%  '
